# Facilitators guidelines - Conditional Probability (Swahili -> English labels,
# plus default-language retag sw-KE -> sw-TZ).

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "Kichwa cha Video" "Video Title"
Replace-Text "Mada" "Topic"
Replace-Text "Malengo" "Aim(s)"
Replace-Text "Urefu" "Length"
Replace-Text "Mahali pa Kambi" "Camp Location"
Replace-Text "Wawezeshaji" "Facilitators"
Replace-Text "N. ya wanafunzi" "N. of students"
Replace-Text "Tarehe" "Date"
Replace-Text "Rasilimali" "Resources"
Replace-Text "inahitajika" "needed"
Replace-Text "Maandalizi" "Preparations"
Replace-Text "Muda wa video" "Video time"
Replace-Text "Mwezeshaji anafanya nini" "What facilitator does"
Replace-Text "Wanachofanya wanafunzi" "What learners do"
Replace-Text "Utangulizi Mkuu wa Video ya VMC" "General VMC Video Introduction"
Replace-Text "Utangulizi wa Video" "Video Introduction"

Replace-Text "Mbona?" "Why?"

# Document default language: Swahili (Kenya) -> Swahili (Tanzania).
$normalStyle = $d.Styles.Item("Normal")
$normalStyle.LanguageID = "sw-TZ"
